$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 4.861952666666666
$ws.Range("H2").Value = 14.585858
$ws.Range("I2").Value = 0.3995648519435639
$ws.Range("J2").Value = 0.3995648519435638
$ws.Range("M2").Value = 0.032838
$ws.Range("N2").Value = 0.098514
$ws.Range("O2").Value = 0.007146324094219707
$ws.Range("P2").Value = 0.007146324094219707
$ws.Range("Q2").Value = 0.159656801668
$ws.Range("R2").Value = 1.436911215012
$ws.Range("S2").Value = 0.00285541992864762
$ws.Range("T2").Value = 0.00285541992864762
$ws.Range("G3").Value = 4.861952666666666
$ws.Range("H3").Value = 14.585858
$ws.Range("I3").Value = 0.3995648519435639
$ws.Range("J3").Value = 0.3995648519435638
$ws.Range("O3").Value = 0.03951718316124263
$ws.Range("P3").Value = 0.03951718316124263
$ws.Range("Q3").Value = 0.8828576749766666
$ws.Range("R3").Value = 7.945719074789999
$ws.Range("S3").Value = 0.0157896774390486
$ws.Range("T3").Value = 0.0157896774390486
$ws.Range("G4").Value = 4.861952666666666
$ws.Range("H4").Value = 14.585858
$ws.Range("I4").Value = 0.3995648519435639
$ws.Range("J4").Value = 0.3995648519435638
$ws.Range("M4").Value = 3.814633
$ws.Range("N4").Value = 11.443899
$ws.Range("O4").Value = 0.8301542030119253
$ws.Range("P4").Value = 0.8301542030119253
$ws.Range("Q4").Value = 18.54656508670467
$ws.Range("R4").Value = 166.919085780342
$ws.Range("S4").Value = 0.3317004412167872
$ws.Range("T4").Value = 0.3317004412167872
$ws.Range("G5").Value = 4.861952666666666
$ws.Range("H5").Value = 14.585858
$ws.Range("I5").Value = 0.3995648519435639
$ws.Range("J5").Value = 0.3995648519435638
$ws.Range("M5").Value = 0.5660336666666667
$ws.Range("N5").Value = 1.698101
$ws.Range("O5").Value = 0.1231822897326124
$ws.Range("P5").Value = 0.1231822897326124
$ws.Range("Q5").Value = 2.752028895073111
$ws.Range("R5").Value = 24.76826005565799
$ws.Range("S5").Value = 0.04921931335908047
$ws.Range("T5").Value = 0.04921931335908045
$ws.Range("I6").Value = 0.04932556406896855
$ws.Range("J6").Value = 0.04932556406896854
$ws.Range("M6").Value = 0.032838
$ws.Range("N6").Value = 0.098514
$ws.Range("O6").Value = 0.007146324094219707
$ws.Range("P6").Value = 0.007146324094219707
$ws.Range("Q6").Value = 0.019709345708
$ws.Range("R6").Value = 0.177384111372
$ws.Range("S6").Value = 0.0003524964669670478
$ws.Range("T6").Value = 0.0003524964669670478
$ws.Range("I7").Value = 0.04932556406896855
$ws.Range("J7").Value = 0.04932556406896854
$ws.Range("O7").Value = 0.03951718316124263
$ws.Range("P7").Value = 0.03951718316124263
$ws.Range("S7").Value = 0.001949207349845038
$ws.Range("T7").Value = 0.001949207349845038
$ws.Range("I8").Value = 0.04932556406896855
$ws.Range("J8").Value = 0.04932556406896854
$ws.Range("M8").Value = 3.814633
$ws.Range("N8").Value = 11.443899
$ws.Range("O8").Value = 0.8301542030119253
$ws.Range("P8").Value = 0.8301542030119253
$ws.Range("Q8").Value = 2.289540183511334
$ws.Range("R8").Value = 20.605861651602
$ws.Range("S8").Value = 0.04094782432778825
$ws.Range("T8").Value = 0.04094782432778824
$ws.Range("I9").Value = 0.04932556406896855
$ws.Range("J9").Value = 0.04932556406896854
$ws.Range("M9").Value = 0.5660336666666667
$ws.Range("N9").Value = 1.698101
$ws.Range("O9").Value = 0.1231822897326124
$ws.Range("P9").Value = 0.1231822897326124
$ws.Range("Q9").Value = 0.3397330293775556
$ws.Range("R9").Value = 3.057597264398
$ws.Range("S9").Value = 0.00607603592436822
$ws.Range("T9").Value = 0.006076035924368218
$ws.Range("G10").Value = 4.206754333333333
$ws.Range("H10").Value = 12.620263
$ws.Range("I10").Value = 0.3457193616641432
$ws.Range("J10").Value = 0.3457193616641432
$ws.Range("M10").Value = 0.032838
$ws.Range("N10").Value = 0.098514
$ws.Range("O10").Value = 0.007146324094219707
$ws.Range("P10").Value = 0.007146324094219707
$ws.Range("Q10").Value = 0.138141398798
$ws.Range("R10").Value = 1.243272589182
$ws.Range("S10").Value = 0.002470622604098724
$ws.Range("T10").Value = 0.002470622604098724
$ws.Range("G11").Value = 4.206754333333333
$ws.Range("H11").Value = 12.620263
$ws.Range("I11").Value = 0.3457193616641432
$ws.Range("J11").Value = 0.3457193616641432
$ws.Range("O11").Value = 0.03951718316124263
$ws.Range("P11").Value = 0.03951718316124263
$ws.Range("Q11").Value = 0.7638834856183332
$ws.Range("R11").Value = 6.874951370564998
$ws.Range("S11").Value = 0.01366185533726983
$ws.Range("T11").Value = 0.01366185533726983
$ws.Range("G12").Value = 4.206754333333333
$ws.Range("H12").Value = 12.620263
$ws.Range("I12").Value = 0.3457193616641432
$ws.Range("J12").Value = 0.3457193616641432
$ws.Range("M12").Value = 3.814633
$ws.Range("N12").Value = 11.443899
$ws.Range("O12").Value = 0.8301542030119253
$ws.Range("P12").Value = 0.8301542030119253
$ws.Range("Q12").Value = 16.04722390282633
$ws.Range("R12").Value = 144.425015125437
$ws.Range("S12").Value = 0.2870003811480884
$ws.Range("T12").Value = 0.2870003811480884
$ws.Range("G13").Value = 4.206754333333333
$ws.Range("H13").Value = 12.620263
$ws.Range("I13").Value = 0.3457193616641432
$ws.Range("J13").Value = 0.3457193616641432
$ws.Range("M13").Value = 0.5660336666666667
$ws.Range("N13").Value = 1.698101
$ws.Range("O13").Value = 0.1231822897326124
$ws.Range("P13").Value = 0.1231822897326124
$ws.Range("Q13").Value = 2.381164580062555
$ws.Range("R13").Value = 21.430481220563
$ws.Range("S13").Value = 0.0425865025746863
$ws.Range("T13").Value = 0.0425865025746863
$ws.Range("G14").Value = 2.499212666666667
$ws.Range("H14").Value = 7.497638
$ws.Range("I14").Value = 0.2053902223233243
$ws.Range("J14").Value = 0.2053902223233243
$ws.Range("M14").Value = 0.032838
$ws.Range("N14").Value = 0.098514
$ws.Range("O14").Value = 0.007146324094219707
$ws.Range("P14").Value = 0.007146324094219707
$ws.Range("Q14").Value = 0.082069145548
$ws.Range("R14").Value = 0.7386223099320001
$ws.Range("S14").Value = 0.001467785094506315
$ws.Range("T14").Value = 0.001467785094506315
$ws.Range("G15").Value = 2.499212666666667
$ws.Range("H15").Value = 7.497638
$ws.Range("I15").Value = 0.2053902223233243
$ws.Range("J15").Value = 0.2053902223233243
$ws.Range("O15").Value = 0.03951718316124263
$ws.Range("P15").Value = 0.03951718316124263
$ws.Range("Q15").Value = 0.4538195320766667
$ws.Range("R15").Value = 4.08437578869
$ws.Range("S15").Value = 0.008116443035079152
$ws.Range("T15").Value = 0.008116443035079152
$ws.Range("G16").Value = 2.499212666666667
$ws.Range("H16").Value = 7.497638
$ws.Range("I16").Value = 0.2053902223233243
$ws.Range("J16").Value = 0.2053902223233243
$ws.Range("M16").Value = 3.814633
$ws.Range("N16").Value = 11.443899
$ws.Range("O16").Value = 0.8301542030119253
$ws.Range("P16").Value = 0.8301542030119253
$ws.Range("Q16").Value = 9.533579112284668
$ws.Range("R16").Value = 85.802212010562
$ws.Range("S16").Value = 0.1705055563192615
$ws.Range("T16").Value = 0.1705055563192615
$ws.Range("G17").Value = 2.499212666666667
$ws.Range("H17").Value = 7.497638
$ws.Range("I17").Value = 0.2053902223233243
$ws.Range("J17").Value = 0.2053902223233243
$ws.Range("M17").Value = 0.5660336666666667
$ws.Range("N17").Value = 1.698101
$ws.Range("O17").Value = 0.1231822897326124
$ws.Range("P17").Value = 0.1231822897326124
$ws.Range("Q17").Value = 1.414638509493111
$ws.Range("R17").Value = 12.731746585438
$ws.Range("S17").Value = 0.02530043787447741
$ws.Range("T17").Value = 0.02530043787447741
